# Append row 20 to the "Artfynd" sheet (new species observation record),
# matching the columns already used by the preceding data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

$ws.Cells.Item($row, 1).Value = 112174920          # A  Id
$ws.Cells.Item($row, 2).Value = 90021               # B  Taxonsorteringsordning
$ws.Cells.Item($row, 3).Value = "Ovaliderad"        # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value = "LC"                # D  Rödlistade
$ws.Cells.Item($row, 5).Value = 6031                # E  TaxonId
$ws.Cells.Item($row, 6).Value = "Blomkålssvamp"      # F  Artnamn
$ws.Cells.Item($row, 7).Value = "Sparassis crispa"   # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value = "(Wulfen:Fr.) Fr."   # H  Auktor
$ws.Cells.Item($row, 9).Value = "'2"                # I  Antal (stored as text, like the rest of the column)
$ws.Cells.Item($row, 9).Style = "Normal"
$ws.Cells.Item($row, 10).Value = "fruktkroppar"     # J  Enhet
$ws.Cells.Item($row, 16).Value = "Osbro, Sm"         # P  Lokalnamn
$ws.Cells.Item($row, 17).Value = 537160.6425027023   # Q  Ost
$ws.Cells.Item($row, 18).Value = 6332349.318422875   # R  Nord
$ws.Cells.Item($row, 19).Value = 3                   # S  Noggrannhet
$ws.Cells.Item($row, 20).Value = "Kalmar"            # T  Län
$ws.Cells.Item($row, 21).Value = "Högsby"            # U  Kommun
$ws.Cells.Item($row, 22).Value = "Småland"           # V  Provins
$ws.Cells.Item($row, 23).Value = "Fagerhult"         # W  Församling
$ws.Cells.Item($row, 25).Value = "'2023-09-18"       # Y  Startdatum (stored as text)
$ws.Cells.Item($row, 25).Style = "Normal"
$ws.Cells.Item($row, 26).Value = "11:19"             # Z  Starttid
$ws.Cells.Item($row, 27).Value = "'2023-09-18"       # AA Slutdatum (stored as text)
$ws.Cells.Item($row, 27).Style = "Normal"
$ws.Cells.Item($row, 28).Value = "11:19"             # AB Sluttid
$ws.Cells.Item($row, 30).Value = $false              # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false              # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false              # AG Ospontan
$ws.Cells.Item($row, 49).Value = "Anders  Henriksson" # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Anders  Henriksson" # AX Observatörer
